$wb = $excel.ActiveWorkbook

# --- Rename sheets -------------------------------------------------------
$wsCars = $wb.Worksheets.Item("Sheet1")
$wsCars.Name = "Cars"

$wsPT = $wb.Worksheets.Item("Sheet2")
$wsPT.Name = "Public_transport"

# --- Public_transport: clear the old placeholder data --------------------
$wsPT.Range("C1").ClearContents()
$wsPT.Range("B4:C4").ClearContents()

# --- Public_transport: enter the transport-type / emission-factor table --
# Data rows first, header row last (matches the order the values were
# authored in, so the shared-string table comes out the same way).
$wsPT.Range("A2").Value = "Bus"
$wsPT.Range("B2").NumberFormat = "@"
$wsPT.Range("B2").Value = "0.105"
$wsPT.Range("B2").ClearFormats()

$wsPT.Range("A3").Value = "Railway"
$wsPT.Range("B3").NumberFormat = "@"
$wsPT.Range("B3").Value = "0.041"
$wsPT.Range("B3").ClearFormats()

$wsPT.Range("A1").Value = "Transport_type"
$wsPT.Range("B1").Value = "kg CO2e/km"

# --- Public_transport: cosmetics ------------------------------------------
$wsPT.Columns.Item(1).AutoFit() | Out-Null

$wsPT.PageSetup.PaperSize = 9
$wsPT.PageSetup.Orientation = 1

# --- Selections (match the saved cursor position in each sheet) ----------
$wsCars.Activate() | Out-Null
$wsCars.Range("D5").Select() | Out-Null

$wsPT.Activate() | Out-Null
$wsPT.Range("F22").Select() | Out-Null
